$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New quantity column (G) for existing rows ---
$ws.Range("G5").Value = 1
$ws.Range("G6").Value = 4

# --- New BOM rows: Nut, lock nut, washer ---
$ws.Range("B25").Value = "Nut"
$ws.Range("D25").Value = "Mcmaster"
$ws.Range("C25").Value = "Nut to replace Brass inserts"
$ws.Range("E25").Value = "95462A029"
$ws.Range("F25").Value = "https://www.mcmaster.com/95462A029/"
$ws.Range("G25").Value = 1

$ws.Range("B26").Value = "lock nut"
$ws.Range("C26").Value = "to keep bolt in place"
$ws.Range("E26").Value = "91102A750"
$ws.Range("F26").Value = "https://www.mcmaster.com/91102A750/"
$ws.Range("G26").Value = 1

$ws.Range("B27").Value = "washer "
$ws.Range("E27").Value = "92141A029"
$ws.Range("F27").Value = "https://www.mcmaster.com/92141A029/"
$ws.Range("G27").Value = 1

# --- Small note cell re-added on F15 ---
$ws.Range("F15").Value = " "

# --- Section header for the new "parts to buy" list ---
$ws.Range("A24").Value = "What we need to buy for a new connection design"

# --- Turn the McMaster part-number URLs into real hyperlinks ---
# Order matters: it matches the order the workbook's relationship ids were
# issued in, so keep F6/F5/F4 (existing rows) before the new F25/F26/F27
# rows, followed by F3 last.
$ws.Hyperlinks.Add($ws.Range("F6"), "https://www.mcmaster.com/97042A216/") | Out-Null
$ws.Hyperlinks.Add($ws.Range("F5"), "https://www.mcmaster.com/97763A348/") | Out-Null
$ws.Hyperlinks.Add($ws.Range("F4"), "https://www.mcmaster.com/93650A145/") | Out-Null
$ws.Hyperlinks.Add($ws.Range("F25"), "https://www.mcmaster.com/95462A029/") | Out-Null
$ws.Hyperlinks.Add($ws.Range("F26"), "https://www.mcmaster.com/91102A750/") | Out-Null
$ws.Hyperlinks.Add($ws.Range("F27"), "https://www.mcmaster.com/92141A029/") | Out-Null
$ws.Hyperlinks.Add($ws.Range("F3"), "https://www.mcmaster.com/94459A380/") | Out-Null

# --- Match the "pasted from the web" look for the new part-number cells ---
$newPartNumbers = $ws.Range("E25:E27")
$newPartNumbers.Font.Name = "Arial"
$newPartNumbers.Font.Color = 3355443

# --- Move the active selection like the author left it ---
$ws.Range("F33").Select() | Out-Null
